$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.39305831775856
$ws.Range("D2").Value = 2.79157380833403
$ws.Range("E2").Value = 10.8037143894589
$ws.Range("F2").Value = 20.27479986733831
$ws.Range("G2").Value = 3.570104529791232
$ws.Range("I2").Value = 24.25631413728395
$ws.Range("M2").Value = 23.74438286295263
$ws.Range("O2").Value = 17.45812196154642
$ws.Range("B3").Value = 9.919740286067597
$ws.Range("D3").Value = 2.796522906301453
$ws.Range("E3").Value = 10.87149116235124
$ws.Range("F3").Value = 19.94573952827516
$ws.Range("G3").Value = 3.573123286525179
$ws.Range("I3").Value = 23.46679190532137
$ws.Range("M3").Value = 22.55226389130305
$ws.Range("O3").Value = 17.29604510821176
$ws.Range("B4").Value = 9.616176361483015
$ws.Range("D4").Value = 2.799870054114073
$ws.Range("E4").Value = 10.92169931525633
$ws.Range("F4").Value = 19.75054696237834
$ws.Range("G4").Value = 3.575073040860725
$ws.Range("I4").Value = 22.97426093668741
$ws.Range("M4").Value = 21.78536397901359
$ws.Range("O4").Value = 17.20367083179469
$ws.Range("B5").Value = 9.489333585493803
$ws.Range("D5").Value = 2.801311881364192
$ws.Range("E5").Value = 10.94429861525945
$ws.Range("F5").Value = 19.67284609762919
$ws.Range("G5").Value = 3.575891866643749
$ws.Range("I5").Value = 22.77195476996678
$ws.Range("M5").Value = 21.46433815426412
$ws.Range("O5").Value = 17.1678619861861
$ws.Range("B6").Value = 9.468085700998131
$ws.Range("D6").Value = 2.801556005494525
$ws.Range("E6").Value = 10.94817974220948
$ws.Range("F6").Value = 19.66005853508802
$ws.Range("G6").Value = 3.576029301394154
$ws.Range("I6").Value = 22.73827653732438
$ws.Range("M6").Value = 21.41052728750918
$ws.Range("O6").Value = 17.1620277820377
$ws.Range("B7").Value = 9.614478252772937
$ws.Range("D7").Value = 2.799889183586393
$ws.Range("E7").Value = 10.92199546646563
$ws.Range("F7").Value = 19.74949145124859
$ws.Range("G7").Value = 3.57508398538606
$ws.Range("I7").Value = 22.97153853598644
$ws.Range("M7").Value = 21.78106855816605
$ws.Range("O7").Value = 17.20318042638613
$ws.Range("B8").Value = 10.23260954278775
$ws.Range("D8").Value = 2.793216426515251
$ws.Range("E8").Value = 10.82528780094209
$ws.Range("F8").Value = 20.1599904165332
$ws.Range("G8").Value = 3.571125478839591
$ws.Range("I8").Value = 23.98590258248302
$ws.Range("M8").Value = 23.34074166599844
$ws.Range("O8").Value = 17.40078027068309
$ws.Range("B9").Value = 11.33779386101084
$ws.Range("D9").Value = 2.782565488521791
$ws.Range("E9").Value = 10.70470910152742
$ws.Range("F9").Value = 21.01392383354163
$ws.Range("G9").Value = 3.564122328779255
$ws.Range("I9").Value = 25.89927787350263
$ws.Range("M9").Value = 26.11226152305616
$ws.Range("O9").Value = 17.84301050933569
$ws.Range("B10").Value = 12.07992748932232
$ws.Range("D10").Value = 2.776207837783647
$ws.Range("E10").Value = 10.6592902474916
$ws.Range("F10").Value = 21.66372716849058
$ws.Range("G10").Value = 3.559434449724971
$ws.Range("I10").Value = 27.2418217314535
$ws.Range("M10").Value = 27.96315177935889
$ws.Range("O10").Value = 18.19841365965101
$ws.Range("B11").Value = 12.40165608269878
$ws.Range("D11").Value = 2.773630884747887
$ws.Range("E11").Value = 10.6481834152161
$ws.Range("F11").Value = 21.96262686407535
$ws.Range("G11").Value = 3.557399908958468
$ws.Range("I11").Value = 27.83582396146429
$ws.Range("M11").Value = 28.7633751860998
$ws.Range("O11").Value = 18.36603636500406
$ws.Range("B12").Value = 12.52115746413475
$ws.Range("D12").Value = 2.772700112408875
$ws.Range("E12").Value = 10.6453626585197
$ws.Range("F12").Value = 22.07615768547625
$ws.Range("G12").Value = 3.556643480708358
$ws.Range("I12").Value = 28.05813562461224
$ws.Range("M12").Value = 29.06029804810245
$ws.Range("O12").Value = 18.43030538893189
$ws.Range("B13").Value = 12.49552501247449
$ws.Range("D13").Value = 2.772898570689518
$ws.Range("E13").Value = 10.64590840330407
$ws.Range("F13").Value = 22.05169353655843
$ws.Range("G13").Value = 3.556805769381387
$ws.Range("I13").Value = 28.01037675291932
$ws.Range("M13").Value = 28.99662326341378
$ws.Range("O13").Value = 18.41642958862549
$ws.Range("B14").Value = 12.41153447838104
$ws.Range("D14").Value = 2.773553407857927
$ws.Range("E14").Value = 10.6479235344328
$ws.Range("F14").Value = 21.97196095530753
$ws.Range("G14").Value = 3.557337396885897
$ws.Range("I14").Value = 27.85416742975812
$ws.Range("M14").Value = 28.7879260483583
$ws.Range("O14").Value = 18.37130826196824
$ws.Range("B15").Value = 12.35978311133749
$ws.Range("D15").Value = 2.773960375622003
$ws.Range("E15").Value = 10.64933853669573
$ws.Range("F15").Value = 21.92316342440156
$ws.Range("G15").Value = 3.557664855876228
$ws.Range("I15").Value = 27.75813691463068
$ws.Range("M15").Value = 28.65929539749902
$ws.Range("O15").Value = 18.34377169088891
$ws.Range("B16").Value = 12.05857808833512
$ws.Range("D16").Value = 2.77638256765523
$ws.Range("E16").Value = 10.66020938704003
$ws.Range("F16").Value = 21.64424854927297
$ws.Range("G16").Value = 3.559569376554411
$ws.Range("I16").Value = 27.20264641328851
$ws.Range("M16").Value = 27.91000603587178
$ws.Range("O16").Value = 18.18757360522778
$ws.Range("B17").Value = 11.8696964475825
$ws.Range("D17").Value = 2.777949042287577
$ws.Range("E17").Value = 10.66933426618506
$ws.Range("F17").Value = 21.47389479034323
$ws.Range("G17").Value = 3.560762778341883
$ws.Range("I17").Value = 26.85742292912881
$ws.Range("M17").Value = 27.43956864583718
$ws.Range("O17").Value = 18.09323148998904
$ws.Range("B18").Value = 11.75956504802281
$ws.Range("D18").Value = 2.778879717801801
$ws.Range("E18").Value = 10.67548145630042
$ws.Range("F18").Value = 21.37623016259614
$ws.Range("G18").Value = 3.561458420609285
$ws.Range("I18").Value = 26.65729642425046
$ws.Range("M18").Value = 27.1650597359735
$ws.Range("O18").Value = 18.03953206544488
$ws.Range("B19").Value = 11.72202185998843
$ws.Range("D19").Value = 2.779199934174632
$ws.Range("E19").Value = 10.67771677263171
$ws.Range("F19").Value = 21.34322130765839
$ws.Range("G19").Value = 3.561695540800839
$ws.Range("I19").Value = 26.58927564629714
$ws.Range("M19").Value = 27.07144487747818
$ws.Range("O19").Value = 18.0214490132226
$ws.Range("B20").Value = 11.88995797877252
$ws.Range("D20").Value = 2.777779218488864
$ws.Range("E20").Value = 10.66826979535619
$ws.Range("F20").Value = 21.49199724807124
$ws.Range("G20").Value = 3.560634784118261
$ws.Range("I20").Value = 26.89433599629282
$ws.Range("M20").Value = 27.4900546015044
$ws.Range("O20").Value = 18.10321647127925
$ws.Range("B21").Value = 12.4362681150963
$ws.Range("D21").Value = 2.773359845310675
$ws.Range("E21").Value = 10.64729397251965
$ws.Range("F21").Value = 21.9953720357814
$ws.Range("G21").Value = 3.557180865490551
$ws.Range("I21").Value = 27.90012272983645
$ws.Range("M21").Value = 28.84939184420978
$ws.Range("O21").Value = 18.38454043937713
$ws.Range("B22").Value = 12.77971456750474
$ws.Range("D22").Value = 2.770734080805183
$ws.Range("E22").Value = 10.64166210834443
$ws.Range("F22").Value = 22.32630361337703
$ws.Range("G22").Value = 3.555005141686459
$ws.Range("I22").Value = 28.54208322468849
$ws.Range("M22").Value = 29.70218391261389
$ws.Range("O22").Value = 18.57300010463693
$ws.Range("B23").Value = 12.5976686449672
$ws.Range("D23").Value = 2.772111560845787
$ws.Range("E23").Value = 10.64392587995841
$ws.Range("F23").Value = 22.1495432913122
$ws.Range("G23").Value = 3.556158926601863
$ws.Range("I23").Value = 28.20092986520537
$ws.Range("M23").Value = 29.25031918028466
$ws.Range("O23").Value = 18.47201557926857
$ws.Range("B24").Value = 11.88080253614841
$ws.Range("D24").Value = 2.777855902060461
$ws.Range("E24").Value = 10.66874823622954
$ws.Range("F24").Value = 21.4838122633777
$ws.Range("G24").Value = 3.560692620584364
$ws.Range("I24").Value = 26.87765273024862
$ws.Range("M24").Value = 27.46724250238895
$ws.Range("O24").Value = 18.09870057727425
$ws.Range("B25").Value = 11.05079525470416
$ws.Range("D25").Value = 2.785187957680931
$ws.Range("E25").Value = 10.72981661554602
$ws.Range("F25").Value = 20.77843494088766
$ws.Range("G25").Value = 3.565936146117284
$ws.Range("I25").Value = 25.3916401407595
$ws.Range("M25").Value = 25.39442234831322
$ws.Range("O25").Value = 17.71779782767767
